$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Text fix: "RuleName" -> "Rule Name" (cell B8)
# ---------------------------------------------------------------------------
$ws.Range("B8").Value = "Rule Name"

# ---------------------------------------------------------------------------
# 2. Colour scheme (Drools/Guvnor decision-table look) applied to the header
#    rows of the rule table.
#    Colors are plain decimal BGR integers (no RGB() helper available):
#      white            = 16777215   (FFFFFF)
#      dark (1C1C1C)    = 1842204
#      dark bg (333300) = 13107
#      yellow (FFFF00)  = 65535
#      orange fg FF860D = 886527
#      orange bg FF6600 = 26367
#      teal   fg 50938A = 9081680
#      teal   bg 808080 = 8421504
#      blue   fg 2A6099 = 10051626
#      blue   bg 666699 = 10053222
# ---------------------------------------------------------------------------

# Row 1 (A1:B1) - "Ruleset" title bar: bold white text on near-black fill
$rng = $ws.Range("A1:B1")
$rng.Font.Bold = $true
$rng.Font.Color = 16777215
$rng.Interior.Color = 1842204
$rng.Interior.PatternColor = 13107

# Row 2 (A2:B2) - "import" line: plain text on yellow fill
$rng = $ws.Range("A2:B2")
$rng.Interior.Color = 65535
$rng.Interior.PatternColor = 65535
# B2 keeps its pre-existing wrap-text behaviour; re-assert it explicitly since
# the load path does not round-trip wrapText on its own.
$ws.Range("B2").WrapText = $true

# Row 5 (C5:G5) - condition/action header band: white text on orange fill
$rng = $ws.Range("C5:G5")
$rng.Font.Color = 16777215
$rng.Interior.Color = 886527
$rng.Interior.PatternColor = 26367

# Row 5 (H5) - action column header: white text on teal fill
$rng = $ws.Range("H5")
$rng.Font.Color = 16777215
$rng.Interior.Color = 9081680
$rng.Interior.PatternColor = 8421504

# Row 8 (B8) - "Rule Name" label: white text on blue fill
$rng = $ws.Range("B8")
$rng.Font.Color = 16777215
$rng.Interior.Color = 10051626
$rng.Interior.PatternColor = 10053222

# Row 8 (C8:H8) - column captions: bold white text on blue fill
$rng = $ws.Range("C8:H8")
$rng.Font.Bold = $true
$rng.Font.Color = 16777215
$rng.Interior.Color = 10051626
$rng.Interior.PatternColor = 10053222

# ---------------------------------------------------------------------------
# 3. Column B width: 61.92 -> 25.72 characters
# ---------------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 24.85

# ---------------------------------------------------------------------------
# 4. Active-cell selection moves from B4 to C8
# ---------------------------------------------------------------------------
$ws.Range("C8").Select() | Out-Null
